# Apply the "output generated at 456a3b4" update to 上海-漫展信息.xlsx
#
# Sheet order in this workbook:
#   1 = 展览      (Exhibitions)
#   2 = 演出      (Performances)
#   3 = 本地生活  (Local life)
#   4 = 全部类型  (All types - aggregate of the above three)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 - "want to go" counts (column F) bumped, plus one cover image
# swapped out.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(3, 6).Value  = 3469
$ws1.Cells.Item(5, 6).Value  = 8127
$ws1.Cells.Item(8, 6).Value  = 2107
$ws1.Cells.Item(9, 6).Value  = 67
$ws1.Cells.Item(10, 6).Value = 540
$ws1.Cells.Item(12, 6).Value = 510
$ws1.Cells.Item(13, 6).Value = 1069
$ws1.Cells.Item(15, 6).Value = 150
$ws1.Cells.Item(16, 6).Value = 1153
$ws1.Cells.Item(18, 6).Value = 723
$ws1.Cells.Item(20, 6).Value = 9
$ws1.Cells.Item(24, 6).Value = 5675
$ws1.Cells.Item(25, 6).Value = 107
$ws1.Cells.Item(26, 6).Value = 51412
$ws1.Cells.Item(27, 6).Value = 4056
$ws1.Cells.Item(29, 6).Value = 995
$ws1.Cells.Item(30, 6).Value = 764
$ws1.Cells.Item(35, 6).Value = 569
$ws1.Cells.Item(36, 6).Value = 1200
$ws1.Cells.Item(38, 6).Value = 4
$ws1.Cells.Item(40, 6).Value = 1026
$ws1.Cells.Item(42, 6).Value = 161
$ws1.Cells.Item(45, 6).Value = 99
$ws1.Cells.Item(47, 6).Value = 98
$ws1.Cells.Item(49, 6).Value = 2453

$ws1.Cells.Item(44, 9).Value = "//i0.hdslb.com/bfs/openplatform/202409/DQLGW65C1726814328151.jpeg"

# ---------------------------------------------------------------------------
# Sheet 2: 演出 - "want to go" counts bumped
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(18, 6).Value = 156
$ws2.Cells.Item(19, 6).Value = 7309
$ws2.Cells.Item(20, 6).Value = 73

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 - a few "want to go" counts bumped, and a brand-new event
# (星零界) inserted as row 14, pushing the old row 14 (咒术回战 0 cafe,
# which also changed address/count) down to row 15.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(4, 6).Value  = 2207
$ws3.Cells.Item(5, 6).Value  = 1484
$ws3.Cells.Item(7, 6).Value  = 630
$ws3.Cells.Item(9, 6).Value  = 9250
$ws3.Cells.Item(10, 6).Value = 1528

# Insert the new row 14; everything from the old row 14 on slides down to 15.
$ws3.Rows.Item(14).Insert()

# Column A carries a bordered/bold numbering style - clone it from the row
# above instead of re-deriving it by hand.
$ws3.Cells.Item(13, 1).Copy($ws3.Cells.Item(14, 1))
$ws3.Cells.Item(14, 1).Value = 13

# Column B holds plain-text dates ("2024-09-24"); force text so Excel does
# not reinterpret the string as a date serial, then drop back to the
# workbook's default (un-styled) look.
$ws3.Cells.Item(14, 2).NumberFormat = "@"
$ws3.Cells.Item(14, 2).Value = "2024-09-24"
$ws3.Cells.Item(14, 2).Style = "Normal"

$ws3.Cells.Item(14, 3).Value = "上海·星零界·社交游乐·休闲运动·潮玩派对"
$ws3.Cells.Item(14, 4).Value = "长宁路1191号长宁来福士B1 长宁来福士"
$ws3.Cells.Item(14, 5).Value = "2024.09.24 10:00-12.31 22:00"
$ws3.Cells.Item(14, 6).Value = 0
$ws3.Cells.Item(14, 7).Value = 68
$ws3.Cells.Item(14, 8).Value = "https://show.bilibili.com/platform/detail.html?id=92659"
$ws3.Cells.Item(14, 9).Value = "//i0.hdslb.com/bfs/openplatform/202409/PHS8s1lu1726221065737.png"

# Former row 14 (now row 15): renumber, and update its address + count.
$ws3.Cells.Item(15, 1).Value = 14
$ws3.Cells.Item(15, 4).Value = "大悦城 次元波板糖"
$ws3.Cells.Item(15, 6).Value = 69

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 - aggregate view mirrors the same "want to go" bumps as
# sheets 1-3 above (no row insertion needed here; it already lists every
# event once).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(2, 6).Value  = 3469
$ws4.Cells.Item(3, 6).Value  = 2207
$ws4.Cells.Item(4, 6).Value  = 1484
$ws4.Cells.Item(6, 6).Value  = 1528
$ws4.Cells.Item(10, 6).Value = 2107
$ws4.Cells.Item(11, 6).Value = 67
$ws4.Cells.Item(13, 6).Value = 510
$ws4.Cells.Item(14, 6).Value = 1069
$ws4.Cells.Item(18, 6).Value = 150
$ws4.Cells.Item(19, 6).Value = 1153
$ws4.Cells.Item(20, 6).Value = 723
$ws4.Cells.Item(21, 6).Value = 9
$ws4.Cells.Item(24, 6).Value = 107
$ws4.Cells.Item(27, 6).Value = 995
$ws4.Cells.Item(30, 6).Value = 569
$ws4.Cells.Item(35, 6).Value = 7309
$ws4.Cells.Item(37, 6).Value = 73
$ws4.Cells.Item(42, 6).Value = 99
$ws4.Cells.Item(44, 6).Value = 98
$ws4.Cells.Item(47, 6).Value = 2453
